$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 676
$ws.Cells.Item(676, 4).Value = 44890
$ws.Cells.Item(676, 11).Value = "Artic Star"
$ws.Cells.Item(676, 12).Value = "Segunda"
$ws.Cells.Item(676, 13).Value = 180
$ws.Cells.Item(676, 14).Value = 16000
$ws.Cells.Item(676, 15).Value = 16000
$ws.Cells.Item(676, 16).Value = 16000
$ws.Cells.Item(676, 19).Value = 1067

# Row 677
$ws.Cells.Item(677, 4).Value = 44890
$ws.Cells.Item(677, 11).Value = "Early Glo"
$ws.Cells.Item(677, 12).Value = "Primera"
$ws.Cells.Item(677, 14).Value = 12000
$ws.Cells.Item(677, 15).Value = 12000
$ws.Cells.Item(677, 16).Value = 12000
$ws.Cells.Item(677, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(677, 19).Value = 1500
$ws.Cells.Item(677, 20).Value = 8

# Row 678
$ws.Cells.Item(678, 4).Value = 44890
$ws.Cells.Item(678, 11).Value = "Early Glo"
$ws.Cells.Item(678, 12).Value = "Segunda"
$ws.Cells.Item(678, 13).Value = 100
$ws.Cells.Item(678, 14).Value = 10000
$ws.Cells.Item(678, 15).Value = 10000
$ws.Cells.Item(678, 16).Value = 10000
$ws.Cells.Item(678, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(678, 19).Value = 1250
$ws.Cells.Item(678, 20).Value = 8

# Row 679
$ws.Cells.Item(679, 11).Value = "August Red"
$ws.Cells.Item(679, 12).Value = "Especial"
$ws.Cells.Item(679, 13).Value = 70
$ws.Cells.Item(679, 17).Value = "`$/caja 15 kilos empedrada"

# Row 680
$ws.Cells.Item(680, 11).Value = "August Red"
$ws.Cells.Item(680, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(680, 13).Value = 120
$ws.Cells.Item(680, 14).Value = 15000
$ws.Cells.Item(680, 15).Value = 15000
$ws.Cells.Item(680, 16).Value = 15000
$ws.Cells.Item(680, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(680, 19).Value = 1000

# Row 681
$ws.Cells.Item(681, 11).Value = "June Pearl"
$ws.Cells.Item(681, 12).Value = "Especial"
$ws.Cells.Item(681, 13).Value = 150
$ws.Cells.Item(681, 14).Value = 14000
$ws.Cells.Item(681, 15).Value = 14000
$ws.Cells.Item(681, 16).Value = 14000
$ws.Cells.Item(681, 19).Value = 933

# Row 682
$ws.Cells.Item(682, 4).Value = 44225
$ws.Cells.Item(682, 12).Value = "Primera"
$ws.Cells.Item(682, 13).Value = 100
$ws.Cells.Item(682, 14).Value = 13000
$ws.Cells.Item(682, 15).Value = 13000
$ws.Cells.Item(682, 16).Value = 13000
$ws.Cells.Item(682, 19).Value = 867

# Row 683
$ws.Cells.Item(683, 4).Value = 44225
$ws.Cells.Item(683, 11).Value = "Venus"
$ws.Cells.Item(683, 12).Value = "Especial"
$ws.Cells.Item(683, 13).Value = 230

# Row 684
$ws.Cells.Item(684, 4).Value = 44225
$ws.Cells.Item(684, 11).Value = "Venus"
$ws.Cells.Item(684, 13).Value = 170
$ws.Cells.Item(684, 14).Value = 12000
$ws.Cells.Item(684, 15).Value = 12000
$ws.Cells.Item(684, 16).Value = 12000
$ws.Cells.Item(684, 19).Value = 800

# Row 685
$ws.Cells.Item(685, 11).Value = "June Pearl"
$ws.Cells.Item(685, 13).Value = 500
$ws.Cells.Item(685, 14).Value = 12000
$ws.Cells.Item(685, 15).Value = 12000
$ws.Cells.Item(685, 16).Value = 12000
$ws.Cells.Item(685, 19).Value = 800

# Row 686
$ws.Cells.Item(686, 11).Value = "June Pearl"
$ws.Cells.Item(686, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(686, 13).Value = 300
$ws.Cells.Item(686, 14).Value = 14000
$ws.Cells.Item(686, 15).Value = 14000
$ws.Cells.Item(686, 16).Value = 14000
$ws.Cells.Item(686, 19).Value = 933

# Row 687
$ws.Cells.Item(687, 11).Value = "June Pearl"
$ws.Cells.Item(687, 12).Value = "Primera"
$ws.Cells.Item(687, 13).Value = 300

# Row 688
$ws.Cells.Item(688, 4).Value = 44236
$ws.Cells.Item(688, 11).Value = "Venus"
$ws.Cells.Item(688, 13).Value = 400
$ws.Cells.Item(688, 14).Value = 13000
$ws.Cells.Item(688, 15).Value = 13000
$ws.Cells.Item(688, 16).Value = 13000
$ws.Cells.Item(688, 19).Value = 867

# Row 689
$ws.Cells.Item(689, 4).Value = 44236
$ws.Cells.Item(689, 11).Value = "Venus"
$ws.Cells.Item(689, 13).Value = 400

# Row 690
$ws.Cells.Item(690, 4).Value = 44236
$ws.Cells.Item(690, 11).Value = "Venus"
$ws.Cells.Item(690, 12).Value = "Segunda"
$ws.Cells.Item(690, 14).Value = 10000
$ws.Cells.Item(690, 15).Value = 10000
$ws.Cells.Item(690, 16).Value = 10000
$ws.Cells.Item(690, 19).Value = 667

# Row 691
$ws.Cells.Item(691, 11).Value = "Artic Snow"
$ws.Cells.Item(691, 12).Value = "Especial"
$ws.Cells.Item(691, 13).Value = 230

# Row 692
$ws.Cells.Item(692, 1).Value = 5
$ws.Cells.Item(692, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(692, 3).Value = "Maule"
$ws.Cells.Item(692, 4).Value = 44263
$ws.Cells.Item(692, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(692, 5).Value = 7
$ws.Cells.Item(692, 6).Value = "Fruta"
$ws.Cells.Item(692, 7).Value = 100103
$ws.Cells.Item(692, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(692, 9).Value = 100103006
$ws.Cells.Item(692, 10).Value = "Nectarín"
$ws.Cells.Item(692, 11).Value = "Artic Snow"
$ws.Cells.Item(692, 12).Value = "Primera"
$ws.Cells.Item(692, 13).Value = 200
$ws.Cells.Item(692, 14).Value = 11000
$ws.Cells.Item(692, 15).Value = 11000
$ws.Cells.Item(692, 16).Value = 11000
$ws.Cells.Item(692, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(692, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(692, 19).Value = 733
$ws.Cells.Item(692, 20).Value = 15

# Row 693
$ws.Cells.Item(693, 1).Value = 5
$ws.Cells.Item(693, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(693, 3).Value = "Maule"
$ws.Cells.Item(693, 4).Value = 44263
$ws.Cells.Item(693, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(693, 5).Value = 7
$ws.Cells.Item(693, 6).Value = "Fruta"
$ws.Cells.Item(693, 7).Value = 100103
$ws.Cells.Item(693, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(693, 9).Value = 100103006
$ws.Cells.Item(693, 10).Value = "Nectarín"
$ws.Cells.Item(693, 11).Value = "August Red"
$ws.Cells.Item(693, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(693, 13).Value = 200
$ws.Cells.Item(693, 14).Value = 15000
$ws.Cells.Item(693, 15).Value = 15000
$ws.Cells.Item(693, 16).Value = 15000
$ws.Cells.Item(693, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(693, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(693, 19).Value = 1000
$ws.Cells.Item(693, 20).Value = 15

# Row 694
$ws.Cells.Item(694, 1).Value = 5
$ws.Cells.Item(694, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(694, 3).Value = "Maule"
$ws.Cells.Item(694, 4).Value = 44263
$ws.Cells.Item(694, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(694, 5).Value = 7
$ws.Cells.Item(694, 6).Value = "Fruta"
$ws.Cells.Item(694, 7).Value = 100103
$ws.Cells.Item(694, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(694, 9).Value = 100103006
$ws.Cells.Item(694, 10).Value = "Nectarín"
$ws.Cells.Item(694, 11).Value = "August Red"
$ws.Cells.Item(694, 12).Value = "Primera"
$ws.Cells.Item(694, 13).Value = 150
$ws.Cells.Item(694, 14).Value = 12000
$ws.Cells.Item(694, 15).Value = 12000
$ws.Cells.Item(694, 16).Value = 12000
$ws.Cells.Item(694, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(694, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(694, 19).Value = 800
$ws.Cells.Item(694, 20).Value = 15
